# Fall 23 Week 4 inputs — append new matchup rows to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (Player_1, Points_1, Player_2, Points_2) for one matchup.
$newRows = @(
    @(3,0,3,3),
    @(5,0,7,3),
    @(5,3,3,0),
    @(6,0,4,2),
    @(2,0,5,3),
    @(3,0,3,3),
    @(5,1,5,2),
    @(3,1,4,2),
    @(5,3,4,0),
    @(6,1,6,2),
    @(4,0,5,3),
    @(3,0,3,3),
    @(5,2,6,0),
    @(6,0,5,2),
    @(2,0,3,3),
    @(3,0,5,3),
    @(6,3,5,0),
    @(3,1,4,2),
    @(3,0,4,3),
    @(6,0,5,2),
    @(3,1,2,2),
    @(5,2,6,0),
    @(3,1,4,2),
    @(7,1,6,2),
    @(5,2,4,1)
)

$startRow = 2250
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

$nextRow = $startRow + $newRows.Count
$ws.Range("A$nextRow").Select()
